# Append the placeholder legend section (page break + 17 explanatory
# paragraphs) right after the final (empty) paragraph of the document body,
# i.e. immediately before the closing <w:sectPr>.
$d = $word.ActiveDocument

$insertPoint = $d.Range($d.Content.End, $d.Content.End)

$xmlPackage = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:br w:type="page"/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>date</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - šiandienos data (YYYY-MM-DD) </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>place</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - vieta ("Utena") </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>client_name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - kliento pavadinimas/vardas pavardė </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>client_code</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - kliento įmonės/asmens kodas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>address</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - adresas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>phone</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - telefonas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>email</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - el. paštas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>user_code</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - aktyvaus naudotojo kodas (įgaliotinio jei juridinis) </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>client_user_code</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - kliento naudotojo kodas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>representative_user_code</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - įgaliotinio naudotojo kodas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>is_juridinis</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>true</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>false</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>is_fizinis</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>true</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>false</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>juridinis_tick</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>☑</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>☐</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (jei reikia </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>checkbox</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">) </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>fizinis_tick</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>☑</w:t>
      </w:r>
      <w:r>
        <w:t>/</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Segoe UI Symbol" w:hAnsi="Segoe UI Symbol" w:cs="Segoe UI Symbol"/>
        </w:rPr>
        <w:t>☐</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (jei reikia </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>checkbox</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">) </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>representative_name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - įgaliotinio vardas pavardė </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>representative_code</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">} - įgaliotinio asmens kodas </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>{</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>representative_position</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>} - įgaliotinio pareigos</w:t>
      </w:r>
      <w:r>
        <w:t>`</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xmlPackage)
